$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2563003333333334
$ws.Range("H2").Value = 0.7689010000000001
$ws.Range("I2").Value = 0.02986826554325775
$ws.Range("J2").Value = 0.02986826554325775
$ws.Range("M2").Value = 31.58970933333333
$ws.Range("N2").Value = 94.76912799999999
$ws.Range("O2").Value = 0.1653281179348216
$ws.Range("P2").Value = 0.1653281179348216
$ws.Range("Q2").Value = 8.096453032036445
$ws.Range("R2").Value = 72.868077288328
$ws.Range("S2").Value = 0.004938064128244286
$ws.Range("T2").Value = 0.004938064128244285
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2563003333333334
$ws.Range("H3").Value = 0.7689010000000001
$ws.Range("I3").Value = 0.02986826554325775
$ws.Range("J3").Value = 0.02986826554325775
$ws.Range("O3").Value = 0.04759211581068901
$ws.Range("P3").Value = 0.04759211581068902
$ws.Range("Q3").Value = 2.330682373753222
$ws.Range("R3").Value = 20.976141363779
$ws.Range("S3").Value = 0.001421493952799135
$ws.Range("T3").Value = 0.001421493952799135
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2563003333333334
$ws.Range("H4").Value = 0.7689010000000001
$ws.Range("I4").Value = 0.02986826554325775
$ws.Range("J4").Value = 0.02986826554325775
$ws.Range("M4").Value = 69.99258933333333
$ws.Range("N4").Value = 209.977768
$ws.Range("O4").Value = 0.3663136922774535
$ws.Range("P4").Value = 0.3663136922774536
$ws.Range("Q4").Value = 17.93912397699645
$ws.Range("R4").Value = 161.452115792968
$ws.Range("S4").Value = 0.01094115463307419
$ws.Range("T4").Value = 0.01094115463307419
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.2563003333333334
$ws.Range("H5").Value = 0.7689010000000001
$ws.Range("I5").Value = 0.02986826554325775
$ws.Range("J5").Value = 0.02986826554325775
$ws.Range("M5").Value = 5.643668666666667
$ws.Range("N5").Value = 16.931006
$ws.Range("O5").Value = 0.02953674277474804
$ws.Range("P5").Value = 0.02953674277474804
$ws.Range("Q5").Value = 1.446474160489556
$ws.Range("R5").Value = 13.018267444406
$ws.Range("S5").Value = 0.0008822112764790742
$ws.Range("T5").Value = 0.0008822112764790742
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.2563003333333334
$ws.Range("H6").Value = 0.7689010000000001
$ws.Range("I6").Value = 0.02986826554325775
$ws.Range("J6").Value = 0.02986826554325775
$ws.Range("M6").Value = 12.31141666666667
$ws.Range("N6").Value = 36.93425
$ws.Range("O6").Value = 0.06443311412377019
$ws.Range("P6").Value = 0.06443311412377019
$ws.Range("Q6").Value = 3.155420195472222
$ws.Range("R6").Value = 28.39878175925
$ws.Range("S6").Value = 0.001924505362427799
$ws.Range("T6").Value = 0.001924505362427799
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.2563003333333334
$ws.Range("H7").Value = 0.7689010000000001
$ws.Range("I7").Value = 0.02986826554325775
$ws.Range("J7").Value = 0.02986826554325775
$ws.Range("M7").Value = 62.44187400000001
$ws.Range("N7").Value = 187.325622
$ws.Range("O7").Value = 0.3267962170785175
$ws.Range("P7").Value = 0.3267962170785176
$ws.Range("Q7").Value = 16.003873120158
$ws.Range("R7").Value = 144.034858081422
$ws.Range("S7").Value = 0.009760836190233265
$ws.Range("T7").Value = 0.009760836190233265
$ws.Range("G8").Value = 6.495645000000001
$ws.Range("I8").Value = 0.7569777503270297
$ws.Range("J8").Value = 0.7569777503270296
$ws.Range("M8").Value = 31.58970933333333
$ws.Range("N8").Value = 94.76912799999999
$ws.Range("O8").Value = 0.1653281179348216
$ws.Range("P8").Value = 0.1653281179348216
$ws.Range("Q8").Value = 205.19553748252
$ws.Range("R8").Value = 1846.75983734268
$ws.Range("S8").Value = 0.1251497067801031
$ws.Range("T8").Value = 0.1251497067801031
$ws.Range("G9").Value = 6.495645000000001
$ws.Range("I9").Value = 0.7569777503270297
$ws.Range("J9").Value = 0.7569777503270296
$ws.Range("O9").Value = 0.04759211581068901
$ws.Range("P9").Value = 0.04759211581068902
$ws.Range("Q9").Value = 59.068535380985
$ws.Range("R9").Value = 531.616818428865
$ws.Range("S9").Value = 0.03602617275967884
$ws.Range("T9").Value = 0.03602617275967884
$ws.Range("G10").Value = 6.495645000000001
$ws.Range("I10").Value = 0.7569777503270297
$ws.Range("J10").Value = 0.7569777503270296
$ws.Range("M10").Value = 69.99258933333333
$ws.Range("N10").Value = 209.977768
$ws.Range("O10").Value = 0.3663136922774535
$ws.Range("P10").Value = 0.3663136922774536
$ws.Range("Q10").Value = 454.64701294012
$ws.Range("R10").Value = 4091.823116461081
$ws.Range("S10").Value = 0.2772913146941746
$ws.Range("T10").Value = 0.2772913146941747
$ws.Range("G11").Value = 6.495645000000001
$ws.Range("I11").Value = 0.7569777503270297
$ws.Range("J11").Value = 0.7569777503270296
$ws.Range("M11").Value = 5.643668666666667
$ws.Range("N11").Value = 16.931006
$ws.Range("O11").Value = 0.02953674277474804
$ws.Range("P11").Value = 0.02953674277474804
$ws.Range("Q11").Value = 36.65926815629
$ws.Range("R11").Value = 329.93341340661
$ws.Range("S11").Value = 0.02235865709761692
$ws.Range("T11").Value = 0.02235865709761692
$ws.Range("G12").Value = 6.495645000000001
$ws.Range("I12").Value = 0.7569777503270297
$ws.Range("J12").Value = 0.7569777503270296
$ws.Range("M12").Value = 12.31141666666667
$ws.Range("N12").Value = 36.93425
$ws.Range("O12").Value = 0.06443311412377019
$ws.Range("P12").Value = 0.06443311412377019
$ws.Range("Q12").Value = 79.97059211375
$ws.Range("R12").Value = 719.7353290237501
$ws.Range("S12").Value = 0.04877443377597632
$ws.Range("T12").Value = 0.04877443377597632
$ws.Range("G13").Value = 6.495645000000001
$ws.Range("I13").Value = 0.7569777503270297
$ws.Range("J13").Value = 0.7569777503270296
$ws.Range("M13").Value = 62.44187400000001
$ws.Range("N13").Value = 187.325622
$ws.Range("O13").Value = 0.3267962170785175
$ws.Range("P13").Value = 0.3267962170785176
$ws.Range("Q13").Value = 405.6002466387301
$ws.Range("R13").Value = 3650.402219748571
$ws.Range("S13").Value = 0.2473774652194798
$ws.Range("T13").Value = 0.2473774652194799
$ws.Range("G14").Value = 1.804372666666667
$ws.Range("H14").Value = 5.413118000000001
$ws.Range("I14").Value = 0.2102747243676212
$ws.Range("J14").Value = 0.2102747243676212
$ws.Range("M14").Value = 31.58970933333333
$ws.Range("N14").Value = 94.76912799999999
$ws.Range("O14").Value = 0.1653281179348216
$ws.Range("P14").Value = 0.1653281179348216
$ws.Range("Q14").Value = 56.99960806901156
$ws.Range("R14").Value = 512.9964726211041
$ws.Range("S14").Value = 0.03476432442896218
$ws.Range("T14").Value = 0.03476432442896218
$ws.Range("G15").Value = 1.804372666666667
$ws.Range("H15").Value = 5.413118000000001
$ws.Range("I15").Value = 0.2102747243676212
$ws.Range("J15").Value = 0.2102747243676212
$ws.Range("O15").Value = 0.04759211581068901
$ws.Range("P15").Value = 0.04759211581068902
$ws.Range("Q15").Value = 16.40817050523578
$ws.Range("R15").Value = 147.673534547122
$ws.Range("S15").Value = 0.01000741903416454
$ws.Range("T15").Value = 0.01000741903416454
$ws.Range("G16").Value = 1.804372666666667
$ws.Range("H16").Value = 5.413118000000001
$ws.Range("I16").Value = 0.2102747243676212
$ws.Range("J16").Value = 0.2102747243676212
$ws.Range("M16").Value = 69.99258933333333
$ws.Range("N16").Value = 209.977768
$ws.Range("O16").Value = 0.3663136922774535
$ws.Range("P16").Value = 0.3663136922774536
$ws.Range("Q16").Value = 126.2927150622916
$ws.Range("R16").Value = 1136.634435560624
$ws.Range("S16").Value = 0.07702651067572715
$ws.Range("T16").Value = 0.07702651067572716
$ws.Range("G17").Value = 1.804372666666667
$ws.Range("H17").Value = 5.413118000000001
$ws.Range("I17").Value = 0.2102747243676212
$ws.Range("J17").Value = 0.2102747243676212
$ws.Range("M17").Value = 5.643668666666667
$ws.Range("N17").Value = 16.931006
$ws.Range("O17").Value = 0.02953674277474804
$ws.Range("P17").Value = 0.02953674277474804
$ws.Range("Q17").Value = 10.18328148185645
$ws.Range("R17").Value = 91.64953333670802
$ws.Range("S17").Value = 0.006210830445677472
$ws.Range("T17").Value = 0.006210830445677471
$ws.Range("G18").Value = 1.804372666666667
$ws.Range("H18").Value = 5.413118000000001
$ws.Range("I18").Value = 0.2102747243676212
$ws.Range("J18").Value = 0.2102747243676212
$ws.Range("M18").Value = 12.31141666666667
$ws.Range("N18").Value = 36.93425
$ws.Range("O18").Value = 0.06443311412377019
$ws.Range("P18").Value = 0.06443311412377019
$ws.Range("Q18").Value = 22.21438372127778
$ws.Range("R18").Value = 199.9294534915
$ws.Range("S18").Value = 0.01354865531252326
$ws.Range("T18").Value = 0.01354865531252326
$ws.Range("G19").Value = 1.804372666666667
$ws.Range("H19").Value = 5.413118000000001
$ws.Range("I19").Value = 0.2102747243676212
$ws.Range("J19").Value = 0.2102747243676212
$ws.Range("M19").Value = 62.44187400000001
$ws.Range("N19").Value = 187.325622
$ws.Range("O19").Value = 0.3267962170785175
$ws.Range("P19").Value = 0.3267962170785176
$ws.Range("Q19").Value = 112.668410701044
$ws.Range("R19").Value = 1014.015696309396
$ws.Range("S19").Value = 0.06871698447056658
$ws.Range("T19").Value = 0.06871698447056658
$ws.Range("I20").Value = 0.002879259762091359
$ws.Range("J20").Value = 0.002879259762091358
$ws.Range("M20").Value = 31.58970933333333
$ws.Range("N20").Value = 94.76912799999999
$ws.Range("O20").Value = 0.1653281179348216
$ws.Range("P20").Value = 0.1653281179348216
$ws.Range("Q20").Value = 0.7804869484986667
$ws.Range("R20").Value = 7.024382536488
$ws.Range("S20").Value = 0.0004760225975120266
$ws.Range("T20").Value = 0.0004760225975120264
$ws.Range("I21").Value = 0.002879259762091359
$ws.Range("J21").Value = 0.002879259762091358
$ws.Range("O21").Value = 0.04759211581068901
$ws.Range("P21").Value = 0.04759211581068902
$ws.Range("S21").Value = 0.0001370300640465088
$ws.Range("T21").Value = 0.0001370300640465088
$ws.Range("I22").Value = 0.002879259762091359
$ws.Range("J22").Value = 0.002879259762091358
$ws.Range("M22").Value = 69.99258933333333
$ws.Range("N22").Value = 209.977768
$ws.Range("O22").Value = 0.3663136922774535
$ws.Range("P22").Value = 0.3663136922774536
$ws.Range("Q22").Value = 1.729306904658667
$ws.Range("R22").Value = 15.563762141928
$ws.Range("S22").Value = 0.001054712274477588
$ws.Range("T22").Value = 0.001054712274477588
$ws.Range("I23").Value = 0.002879259762091359
$ws.Range("J23").Value = 0.002879259762091358
$ws.Range("M23").Value = 5.643668666666667
$ws.Range("N23").Value = 16.931006
$ws.Range("O23").Value = 0.02953674277474804
$ws.Range("P23").Value = 0.02953674277474804
$ws.Range("Q23").Value = 0.1394381217473334
$ws.Range("R23").Value = 1.254943095726
$ws.Range("S23").Value = 0.0000850439549745747
$ws.Range("T23").Value = 0.0000850439549745747
$ws.Range("I24").Value = 0.002879259762091359
$ws.Range("J24").Value = 0.002879259762091358
$ws.Range("M24").Value = 12.31141666666667
$ws.Range("N24").Value = 36.93425
$ws.Range("O24").Value = 0.06443311412377019
$ws.Range("P24").Value = 0.06443311412377019
$ws.Range("Q24").Value = 0.3041781715833334
$ws.Range("R24").Value = 2.73760354425
$ws.Range("S24").Value = 0.0001855196728428119
$ws.Range("T24").Value = 0.0001855196728428119
$ws.Range("I25").Value = 0.002879259762091359
$ws.Range("J25").Value = 0.002879259762091358
$ws.Range("M25").Value = 62.44187400000001
$ws.Range("N25").Value = 187.325622
$ws.Range("O25").Value = 0.3267962170785175
$ws.Range("P25").Value = 0.3267962170785176
$ws.Range("Q25").Value = 1.542751380918
$ws.Range("R25").Value = 13.884762428262
$ws.Range("S25").Value = 0.0009409311982378483
$ws.Range("T25").Value = 0.0009409311982378484
